# mise en commentaires des paramètres de recherche qu'on n'utilise plus
# - Bump the "Date" metadata value to the new export timestamp.
# - Insert a new "Jurisdiction" metadata row (empty value) right after
#   "Contact" and before "Description" on the Metadata sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# 1) Update the Date property's value (row 8, column B).
$ws.Range("B8").Value = "2024-07-01T07:50:29+00:00"

# 2) Insert a new row above the current "Description" row (row 11) so the
#    new "Jurisdiction" property lands between "Contact" (row 10) and
#    "Description" (which shifts from row 11 down to row 12).
$ws.Rows("11:11").Insert()

# Copy the formatting of the row that is now below the freshly inserted
# blank row (the "Description" row) so the new row matches the rest of
# the table's look (borders/alignment) instead of Excel's bare default.
$ws.Range("A12:B12").Copy()
$ws.Range("A11:B11").PasteSpecial(-4122)

# Fill in the new property/value pair.
$ws.Range("A11").Value = "Jurisdiction"
$ws.Range("B11").Value = ""
